$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("ProductLoanInput")

# Update the repayment strategy value on ProductLoanInput!B17
$r = $ws1.Range("B17")
$r.Value = "Penalties, Fees, Interest, Principal order"

# Give it a new left/top aligned style (this produces a new cellXfs entry)
$r.HorizontalAlignment = -4131
$r.VerticalAlignment = -4160

# Move the active sheet/selection to ProductLoanInput!B17
$ws1.Activate()
$ws1.Range("B17").Select()
